$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "46.624.99"
$ws.Range("E2").Value = "  +6.31%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "2.291.49"
$ws.Range("E3").Value = "  +3.10%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 (BNB)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.02%  "

# Row 6 (Solana)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +12.01%  "

# Row 7 (XRP)
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.566"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.61%  "

# Row 8 (USDC)
$ws.Range("E8").Value = "  +0.11%  "

# Row 9 (Cardano)
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.519"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.86%  "

# Row 10 (Avalanche)
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.40%  "

# Row 11 (Dogecoin)
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0791"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.02%  "

# Row 12 (Polkadot)
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.16%  "

# Row 13 (TRON)
$ws.Range("E13").Value = "  +0.21%  "

# Row 14 (WrappedliquidstakedEther2.0)
$ws.Range("D14").Value = "2.644.10"
$ws.Range("E14").Value = "  +3.16%  "

# Row 15 (WrappedEther)
$ws.Range("D15").Value = "2.293.79"
$ws.Range("E15").Value = "  +3.05%  "

# Row 16 (Chainlink)
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.54%  "

# Row 17 (Polygon)
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.811"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.78%  "

# Row 18 (WrappedBTC)
$ws.Range("D18").Value = "46.636.26"
$ws.Range("E18").Value = "  +6.59%  "

# Row 19 (InternetComputer(DFINITY))
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +13.83%  "

# Row 20 (ShibaInu)
$ws.Range("D20").Value = "0.0₃0938"
$ws.Range("E20").Value = "  +4.15%  "

# Row 21 (Uniswap)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.87%  "

# Row 22 (Litecoin)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.95%  "

# Row 23 (BitcoinCash)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "247.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.06%  "

# Row 24 (PancakeSwap)
$ws.Range("E24").Value = "  +4.10%  "

# Row 25 (Dai)
$ws.Range("E25").Value = "  -0.14%  "

# Row 26 (ImmutableX)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.57%  "

# Row 27 (InjectiveProtocol)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "43.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.16%  "

# Row 28 (Toncoin)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.25%  "

# Row 29 (Cosmos)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.49%  "

# Row 30 (EthereumClassic)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.12%  "

# Row 31 (WEMIXToken)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.55%  "

# Row 32 (Filecoin)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.61%  "

# Row 33 (Monero)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "147.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.81%  "

# Row 34 (Hedera)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0794"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.30%  "

# Row 35 (LidoDAOToken)
$ws.Range("E35").Value = "  +13.79%  "

# Row 36 (Kaspa)
$ws.Range("E36").Value = "  +12.49%  "

# Row 37 (Stellar)
$ws.Range("E37").Value = "  +1.06%  "

# Row 38 (ARBITRUM)
$ws.Range("E38").Value = "  +5.85%  "

# Row 39 (Celestia)
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +22.09%  "

# Row 40 (RenderToken)
$ws.Range("E40").Value = "  +12.07%  "

# Row 41 (NEARProtocol)
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.75%  "

# Row 42 (VeChain)
$ws.Range("E42").Value = "  +1.74%  "

# Row 43 (FirstDigitalUSD)
$ws.Range("E43").Value = "  -0.10%  "

# Row 44 (Stacks)
$ws.Range("E44").Value = "  +11.96%  "

# Row 45 (Maker)
$ws.Range("D45").Value = "1.840.01"
$ws.Range("E45").Value = "  +0.60%  "

# Row 46 (BitcoinSV)
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "87.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +19.91%  "

# Row 47 (Algorand)
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.194"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.53%  "

# Row 48 (ordi)
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "73.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.38%  "

# Row 49 (THORChain)
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.81%  "

# Row 50 (Aave)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "95.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.44%  "

# Row 51 (MultiversX)
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.70%  "
